$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.565.65'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '3.782.93'
$ws.Range('E3').Value = '  -1.51%  '
$ws.Range('E4').Value = '  -0.35%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '646.24'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.50%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '166.27'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.35%  '
$ws.Range('D7').Value = '3.783.19'
$ws.Range('E8').Value = '  -0.03%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.527'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +0.88%  '
$ws.Range('E10').Value = '  -2.36%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.457'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('E12').Value = '  +3.78%  '
$ws.Range('E13').Value = '  -4.91%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '35.02'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -3.05%  '
$ws.Range('D15').Value = '4.418.43'
$ws.Range('E15').Value = '  -1.42%  '
$ws.Range('D16').Value = '3.754.40'
$ws.Range('E16').Value = '  -3.64%  '
$ws.Range('D17').Value = '69.449.00'
$ws.Range('E17').Value = '  -0.62%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '17.79'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -2.12%  '
$ws.Range('E19').Value = '  +0.11%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.03'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.81%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '469.51'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.15%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '9.61'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.99%  '
$ws.Range('E23').Value = '  -0.25%  '
$ws.Range('E24').Value = '  -5.43%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '81.99'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -2.29%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '12.33'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +2.41%  '
$ws.Range('E27').Value = '  +2.84%  '
$ws.Range('E28').Value = '  -3.42%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').Value = '3.930.64'
$ws.Range('E30').Value = '  -1.43%  '
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('E32').Value = '  +2.18%  '
$ws.Range('E33').Value = '  -2.36%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '28.79'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -2.20%  '
$ws.Range('E35').Value = '  +15.00%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').Value = '3.736.94'
$ws.Range('E37').Value = '  -1.21%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '8.89'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -2.28%  '
$ws.Range('E39').Value = '  -2.55%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '3.30'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -4.48%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '5.85'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -1.57%  '
$ws.Range('E42').Value = '  -0.06%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.959'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -2.65%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '45.38'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +3.23%  '
$ws.Range('E46').Value = '  +2.92%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '157.32'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.14%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '47.70'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.57%  '
$ws.Range('E49').Value = '  -1.52%  '
$ws.Range('E50').Value = '  -0.98%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '8.40'
$c.Style = 'Normal'
